$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data values scraped on Tue Feb  6 10:49:59 UTC 2024

$ws.Range("D2").Value = "43.038.85"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.325.04"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.82"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.04"
$ws.Range("E6").Value = "  -3.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.98"
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.82"
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0782"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "2.699.29"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "2.318.37"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "42.757.33"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.04"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "0.0₃0886"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.93"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.00"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("E25").Value = "  +7.36%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.56"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  -14.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.09"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.47"
$ws.Range("E30").Value = "  -5.07%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "139.78"
$ws.Range("E32").Value = "  -15.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.98"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.62"
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0697"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.37"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.29"
$ws.Range("E40").Value = "  +23.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.108"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "1.937.52"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0279"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.21"
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").Value = "2.562.29"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.66"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.03"
$ws.Range("E51").Value = "  -0.02%  "
